$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 10

# Delete rows 3 through 6 (entire rows) so remaining data is only header + row 2
$ws.Range("A3:C6").EntireRow.Delete()
